$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 60
$wsExpo.Range("F4").Value = 1479
$wsExpo.Range("F5").Value = 535
$wsExpo.Range("F7").Value = 10976
$wsExpo.Range("F8").Value = 10976
$wsExpo.Range("F12").Value = 1059
$wsExpo.Range("F13").Value = 746
$wsExpo.Range("F14").Value = 12196
$wsExpo.Range("F15").Value = 12701
$wsExpo.Range("F22").Value = 18

# Sheet "全部类型" (sheet4 / rId4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 60
$wsAll.Range("F5").Value = 1479
$wsAll.Range("F6").Value = 535
$wsAll.Range("F8").Value = 10976
$wsAll.Range("F9").Value = 10976
$wsAll.Range("F13").Value = 1059
$wsAll.Range("F14").Value = 746
$wsAll.Range("F15").Value = 12196
$wsAll.Range("F16").Value = 12701
$wsAll.Range("F23").Value = 18
